# Generate Report for Handoff
# Updates the handoff GUID/hash identifiers and their associated
# timestamps across the Overview / zh-cn / de-de sheets, keeping the
# displayed hyperlink text in sync with the new file names.

$wb = $excel.ActiveWorkbook

$oldGuid = "5ca52f3e-b3ff-414c-8c66-6a3215c46356"
$newGuid = "40509138-7ff3-4f5c-be48-38e1b1058a4a"

$oldHash = "a3ef48961fa209730c8a1e4f29fa49c0d713ac71"
$newHash = "1e7643736535e32690f6842d8c3fd40902d31892"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = "2016-00-13 07:00:32"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    }
}

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$oldZh = "$oldGuid.$oldHash.zh-cn.xlf"
$newZh = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newZh
$ws.Range("E2").Value = "2016-03-13 07:00:29"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldZh) {
        $hl.TextToDisplay = $newZh
    }
}

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$oldDe = "$oldGuid.$oldHash.de-de.xlf"
$newDe = "$newGuid.$newHash.de-de.xlf"
$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newDe
$ws.Range("E2").Value = "2016-03-13 07:00:32"
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldDe) {
        $hl.TextToDisplay = $newDe
    }
}
